$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data: the uploaded file rearranges the table -----------------------
# Row 1 (header) stays "Package" / "User"
# Row 2/3 now list the new package names against users A and B
$ws.Range("A1").Value = "Package"
$ws.Range("B1").Value = "User"
$ws.Range("A2").Value = "openxlsx==4.2.4"
$ws.Range("B2").Value = "A"
$ws.Range("A3").Value = "tidyverse"
$ws.Range("B3").Value = "B"

# --- Formatting: the two package-name cells get an explicit Calibri font ---
$pkgRange = $ws.Range("A2:A3")
$pkgRange.Font.Size = 11
$pkgRange.Font.Color = 0
$pkgRange.Font.Name = "Calibri"

# --- View state: selection moved to H9 -----------------------------------
$ws.Range("H9").Select() | Out-Null
